# Raw and Clean Data from SSA for July 16th
# Adds the 2020-07-16 (serial 44028) row to the daily tracking sheets and
# fills in the corresponding "AU" column on control_obs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# out_vars: append row 47 (copy formatting from row 46, then set values)
# ---------------------------------------------------------------------
$wsOut = $wb.Worksheets.Item("out_vars")
$wsOut.Range("A46:J46").Copy()
$wsOut.Range("A47:J47").PasteSpecial(-4122)
$outVals = @(44028, 324041, 375455, 82567, 37574, 28.953743507765996, 93822, 7801, 9092, 782063)
for ($i = 0; $i -lt $outVals.Length; $i++) {
    $wsOut.Cells.Item(47, $i + 1).Value2 = $outVals[$i]
}

# ---------------------------------------------------------------------
# dates_dx: append row 47
# ---------------------------------------------------------------------
$wsDx = $wb.Worksheets.Item("dates_dx")
$wsDx.Range("A46:L46").Copy()
$wsDx.Range("A47:L47").PasteSpecial(-4122)
$dxVals = @(44028, 0, 1, 0, 0, 1, 0, 0, 0, 1, 0, 4)
for ($i = 0; $i -lt $dxVals.Length; $i++) {
    $wsDx.Cells.Item(47, $i + 1).Value2 = $dxVals[$i]
}

# ---------------------------------------------------------------------
# dates_sx: append row 47
# ---------------------------------------------------------------------
$wsSx = $wb.Worksheets.Item("dates_sx")
$wsSx.Range("A46:N46").Copy()
$wsSx.Range("A47:N47").PasteSpecial(-4122)
$sxVals = @(44028, 0, 1, 0, 0, 1, 0, 1, 0, 0, 1, 1, 0, 0)
for ($i = 0; $i -lt $sxVals.Length; $i++) {
    $wsSx.Cells.Item(47, $i + 1).Value2 = $sxVals[$i]
}

# ---------------------------------------------------------------------
# dates_deaths: row 47 already exists (blank placeholders) - fill it in
# ---------------------------------------------------------------------
$wsDeaths = $wb.Worksheets.Item("dates_deaths")
$wsDeaths.Range("A46").Copy()
$wsDeaths.Range("A47").PasteSpecial(-4122)
$deathVals = @(44028, 0, 0, 0, 0, 2, 1, 1, 1, 2)
for ($i = 0; $i -lt $deathVals.Length; $i++) {
    $wsDeaths.Cells.Item(47, $i + 1).Value2 = $deathVals[$i]
}

# ---------------------------------------------------------------------
# control_obs: fill in the new "AU" column (col 47) for 2020-07-16
# ---------------------------------------------------------------------
$wsCtrl = $wb.Worksheets.Item("control_obs")

# AU1 header date - copy AT1's date format, then write the new serial date
$wsCtrl.Range("AT1").Copy()
$wsCtrl.Range("AU1").PasteSpecial(-4122)
$wsCtrl.Cells.Item(1, 47).Value2 = 44028

# AU9 / AU17 are the blank "separator" rows - copy their special fill from AT
$wsCtrl.Range("AT9").Copy()
$wsCtrl.Range("AU9").PasteSpecial(-4122)
$wsCtrl.Range("AT17").Copy()
$wsCtrl.Range("AU17").PasteSpecial(-4122)

# AU20 total - copy AT20's format, then put the SUM formula in
$wsCtrl.Range("AT20").Copy()
$wsCtrl.Range("AU20").PasteSpecial(-4122)
$wsCtrl.Range("AU20").Formula = "=SUM(AU2:AU18)"

# Regular daily counts (format already matches column default, s=48)
$ctrlVals = @{
    2 = 4406; 3 = 4213; 4 = 4213; 5 = 4213; 6 = 4213; 7 = 3495; 8 = 6126;
    10 = 186; 11 = 186; 12 = 186; 13 = 186; 14 = 186; 15 = 121; 16 = 198;
    18 = 1031
}
foreach ($r in $ctrlVals.Keys) {
    $wsCtrl.Cells.Item($r, 47).Value2 = $ctrlVals[$r]
}

# ---------------------------------------------------------------------
# Restore selections on each sheet, ending on control_obs (keeps it the
# active tab, matching the saved workbook view).
# ---------------------------------------------------------------------
$wsOut.Activate()
$wsOut.Range("C52").Select()

$wsDx.Activate()
$wsDx.Range("C56").Select()

$wsSx.Activate()
$wsSx.Range("F59").Select()

$wsDeaths.Activate()
$wsDeaths.Range("L47").Select()

$wsCtrl.Activate()
$wsCtrl.Range("AQ13").Select()

Write-Output "edit complete"
